# EventCardData.xlsx rules rewrite
# The cards' names/effects are updated to reflect the new "slot" based
# combat rules (player zone vs room zone, free ordering of resolution, etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: "出口" (Exit) -> "房间出口" (Room Exit) ---
$ws.Range("A2").Value = "房间出口"

# --- Row 3: "阴影" (Shadow) -> "隐蔽处" (Hidden spot) ---
$ws.Range("A3").Value = "隐蔽处"

# --- Row 3 effect text (monster flip ability) ---
$ws.Range("D3").Value = "选房间区1张怪物牌横置。横置1张“敏捷”牌，可以改为选房间区任意怪物牌横置。"

# --- Row 4: "拐角" (Corner) keeps its name, only the effect text changes ---
$ws.Range("D4").Value = "选场上1张怪物牌移动到其同区域内的一个空槽位中。横置1张“敏捷”牌，可以改为交换场上两个同区域槽位的所有牌。"

# --- Row 5: "杂物堆" (Junk pile) keeps its name, only the effect text changes ---
$ws.Range("D5").Value = "翻开遭遇牌堆前3张牌，获得其中1张战利品牌。横置1张“感知”牌，可以改为翻开遭遇牌堆前5张牌，获得其中所有战利品牌。"

# --- Row 6: "拉杆" (Lever) keeps its name, only the effect text changes ---
$ws.Range("D6").Value = "选场上1张陷阱牌，将其移动到场上任意位置或送墓。横置1张“智力”牌，可以改为将弃牌堆顶端的1张陷阱牌移动到场上任意位置或送墓。"

# --- Row 7: "宝箱" (Chest) keeps its name, only the effect text changes ---
$ws.Range("D7").Value = "从遭遇牌堆翻开5张牌，获得其中的战利品牌。使用1张《敏捷》发动本牌时，可以额外翻开2张牌。"

# --- Row 2 effect text (room exit ability) ---
$ws.Range("D2").Value = "重抽房间区所有牌。横置1张“体质”牌，可以再重抽任意张手牌。"

# --- Row 8: "冒险者尸体" (Dead adventurer) keeps its name, only the effect text changes ---
$ws.Range("D8").Value = "消耗3时间，将主牌堆第1张怪物牌放在房间区任意空槽位，然后获得遭遇牌堆第1张战利品牌，再获得遗物牌堆顶的1张遗物牌。"

# --- Row 9: "商人" (Resource merchant) keeps name and effect text unchanged ---

# --- Row 10: "训练师" (Trainer) keeps its name, only the effect text changes ---
$ws.Range("D10").Value = "消耗2时间，获得1SP。"

# --- Row 11: "祭坛" (Altar) keeps its name, only the effect text changes ---
$ws.Range("D11").Value = "多选：①将1张手牌送墓，翻开遭遇牌堆顶1张牌，如果是战利品牌则可以获得。②弃置1张战利品牌，获得1道具点。③受到3伤害，获得1SP。"

# --- Wrap-text formatting for the long effect cells in rows 2 and 11 ---
$ws.Range("D2").WrapText = $true
$ws.Range("D11").WrapText = $true

# --- Row heights recalculated by Excel for the new wrapped text ---
$ws.Rows(2).RowHeight = 28.5
$ws.Rows(3).RowHeight = 42.75
$ws.Rows(5).RowHeight = 57
$ws.Rows(11).RowHeight = 57

# --- Selection moved on save ---
$ws.Range("D12").Select()
